$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$headers = @("name", "employee_id", "asset", "arrival", "status", "division", "job_title")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.HorizontalAlignment = -4108  # xlCenter
    $cell.VerticalAlignment = -4160    # xlTop
    $cell.Borders.LineStyle = 1        # xlContinuous
    $cell.Borders.Weight = 2           # xlThin
}

# Data row
$values = @("Roger Doddy", "DODROG", "PT-07S", "04:45 AM", "On Time", "TEXDIST", "Select Maintenance Employee")
for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item(2, $i + 1).Value = $values[$i]
}
